# Daily attendance processing - 2025-11-10 06:58:33
# Swaps the order of the two names/emails listed in column G ("Recorded By")
# for rows where the value is exactly "System, dnasr281@gmail.com",
# "System, admin@admin.com", or "dnasr281@gmail.com, admin@admin.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows whose G value is "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
$rowsSystemDnasr = @(3,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,136,137,138,142,144,145,146,148,151,153)

# Rows whose G value is "System, admin@admin.com" -> "admin@admin.com, System"
$rowsSystemAdmin = @(7,33,59)

# Rows whose G value is "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"
$rowsDnasrAdmin = @(87,113,139)

foreach ($r in $rowsSystemDnasr) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

foreach ($r in $rowsSystemAdmin) {
    $ws.Cells.Item($r, 7).Value = "admin@admin.com, System"
}

foreach ($r in $rowsDnasrAdmin) {
    $ws.Cells.Item($r, 7).Value = "admin@admin.com, dnasr281@gmail.com"
}
